# Add a new "isotope" column (E) to the periodic-table style data sheet.
# Existing columns E (etymology) and F (stateOfMatter) shift right to F and G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; this shifts the old E/F columns (and their
# widths/styles) one slot to the right, matching the OOXML diff exactly.
$ws.Columns.Item(5).Insert()

# Header
$ws.Range("E1").Value = "isotope"

# Per-row isotope values
$ws.Range("E2").Value = "195 Au, 198 Au"
$ws.Range("E3").Value = "107 Ag, 109 Ag"
$ws.Range("E4").Value = "190 Pt, 192 pt, 194pt, 195 pt, 196 pt, 198 pt"
$ws.Range("E5").Value = "36 Ar, 38 Ar, 40 Ar"
$ws.Range("E6").Value = "10 B, 11 B"
$ws.Range("E7").Value = "196 Hg, 198 Hg, 199 Hg, 200 Hg, 201 Hg, 202 Hg, 204 Hg"

# Match the new column width used for the wrapped-text columns (E & F).
$ws.Columns.Item(5).ColumnWidth = 52.92

# Update the view so the active cell / selection matches the post-edit state.
$null = $ws.Range("F7").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
